$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row dates (A23, A24) ---
$ws.Cells.Item(23,1).Value2 = 43108
$ws.Cells.Item(24,1).Value2 = 43110

# --- Add new row 25 ---
$ws.Cells.Item(25,1).Value2 = 43116
$ws.Cells.Item(25,2).Value = "Giovanni"
$ws.Cells.Item(25,3).Value = "Compreso ed implementato il pattern MVC. Iniziato a modellare il progetto secondo questo standard."
$ws.Cells.Item(25,4).Value2 = (2/24)

# Match formatting of similar rows (row 16, which also uses column A/B/C/D with wrapped description + ht 72->72.9)
$ws.Cells.Item(25,3).WrapText = $true
$ws.Rows.Item(25).RowHeight = 72.9

# The "Totale ore" total cell (I2) no longer carries the manual underline
$ws.Cells.Item(2,9).Font.Underline = $false

# Move selection to reflect last user action (D26, just below new data)
$ws.Range("D26").Select() | Out-Null
